$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '30.445.46'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.78%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.094.85'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.57%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '329.75'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('E6').Value = '  -0.01%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.5210'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -1.02%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '53.21'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +18.22%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.08857'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('E11').Value = '  -2.05%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '24.41'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -1.72%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.092.83'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -0.84%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '6.679'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -1.93%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '7.679'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('E16').Value = '  -1.84%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('E18').Value = '  -1.58%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.06586'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -0.30%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '19.24'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +0.15%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -0.06%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.267'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -2.39%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '30.479.10'
$c.Style = "Normal"
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '12.25'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +2.16%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.339'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +3.47%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.333.88'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -1.01%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '22.25'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -2.72%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.578'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +1.62%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '162.36'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.61%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '131.63'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -1.55%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.193'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.34%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.1068'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +0.02%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.668'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +9.14%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '6.125'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -1.38%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '3.883'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.69%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '10.08'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +5.75%  '
$ws.Range('E37').Value = '  -0.61%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.06806'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +0.73%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '12.72'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.54%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.449'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -3.14%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.2260'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +0.71%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.6890'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +1.52%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.262'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('E44').Value = '  -0.02%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '13.97'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -1.66%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.6363'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +0.68%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.198'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -2.44%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '3.619'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('E49').Value = '  +8.34%  '
$ws.Range('E50').Value = '  -2.77%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '81.91'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -1.21%  '
